$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 03:22"

# Row 8 - Alemania: solo cambia Casos criticos (F)
$ws.Cells.Item(8, 6).Value = 2570

# Row 14 - Brasil: actualizacion de datos
$ws.Cells.Item(14, 2).Value = 59324
$ws.Cells.Item(14, 3).Value = 128
$ws.Cells.Item(14, 5).Value = 26107
$ws.Cells.Item(14, 7).Value = 12
$ws.Cells.Item(14, 8).Value = 4057

# Row 43 - Noruega: actualizacion de datos
$ws.Cells.Item(43, 2).Value = 7499
$ws.Cells.Item(43, 3).Value = 6
$ws.Cells.Item(43, 5).Value = 7266

# Filas 76-77: Nueva Zelanda ahora supera a Afganistan, se reordenan
$ws.Cells.Item(76, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(76, 2).Value = 1470
$ws.Cells.Item(76, 3).Value = 9
$ws.Cells.Item(76, 4).Value = 1142
$ws.Cells.Item(76, 5).Value = 310
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 18

$ws.Cells.Item(77, 1).Value = "Afganistan"
$ws.Cells.Item(77, 2).Value = 1463
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 188
$ws.Cells.Item(77, 5).Value = 1228
$ws.Cells.Item(77, 6).Value = 7
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 47

# Filas 103-104: Uruguay ahora supera a Honduras, se reordenan
$ws.Cells.Item(103, 1).Value = "Uruguay"
$ws.Cells.Item(103, 2).Value = 596
$ws.Cells.Item(103, 3).Value = 33
$ws.Cells.Item(103, 4).Value = 370
$ws.Cells.Item(103, 5).Value = 212
$ws.Cells.Item(103, 6).Value = 9
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 14

$ws.Cells.Item(104, 1).Value = "Honduras"
$ws.Cells.Item(104, 2).Value = 591
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 58
$ws.Cells.Item(104, 5).Value = 478
$ws.Cells.Item(104, 6).Value = 10
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 55

# Filas 189-191: Republica de Africa Central y Curazao se adelantan a Dominica
$ws.Cells.Item(189, 1).Value = "Republica de Africa Central"

$ws.Cells.Item(190, 1).Value = "Curazao"
$ws.Cells.Item(190, 4).Value = 11
$ws.Cells.Item(190, 5).Value = 4
$ws.Cells.Item(190, 8).Value = 1

$ws.Cells.Item(191, 1).Value = "Dominica"
$ws.Cells.Item(191, 4).Value = 13
$ws.Cells.Item(191, 5).Value = 3
$ws.Cells.Item(191, 8).Value = 0
